# Convert the poi-category controlled vocabulary workbook:
#  - add English (col C) and German (col D) translations
#  - rename headers to codice_1_livello / label_ITA_1_livello / label_ENG_1_livello / label_DEU_1_livello
#  - apply the "new" font style to the German column (rows 4-10)
#  - resize columns and move the active selection, matching the authored commit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows: fill English (C) and German (D) translations first so that
#     new shared strings are appended in the same order as the target file. ---

$ws.Range("C4").Value2 = "Entertainment"
$ws.Range("D4").Value2 = "Unterhaltungssektor"

$ws.Range("C5").Value2 = "Leisure"
$ws.Range("D5").Value2 = "Freizeitssektor"

$ws.Range("C6").Value2 = "Food"
$ws.Range("D6").Value2 = "Lebensmittelsektor"

$ws.Range("C7").Value2 = "Shopping"
$ws.Range("D7").Value2 = "Einkaufssektor"

$ws.Range("C8").Value2 = "Automotive"
$ws.Range("D8").Value2 = "Automobilsektor"

$ws.Range("C9").Value2 = "Travel/Tourism"
$ws.Range("D9").Value2 = "Reise-/Tourismussektor"

$ws.Range("C10").Value2 = "Geographical"
$ws.Range("D10").Value2 = "Geographischer Sektor"

$ws.Range("C11").Value2 = "Other public services"
$ws.Range("D11").Value2 = "Sonstige öffentliche Dienstleistungen"

# --- Header row: overwrite the two existing headers and add the two new ones. ---

$ws.Range("A3").Value2 = "codice_1_livello"
$ws.Range("B3").Value2 = "label_ITA_1_livello"
$ws.Range("C3").Value2 = "label_ENG_1_livello"
$ws.Range("D3").Value2 = "label_DEU_1_livello"

# C3/D3 are brand new cells, so give them the same bold header style already
# used by A3/B3 (this reuses the existing bold font, no new style is created).
$ws.Range("C3:D3").Font.Bold = $true

# --- New font/style used for the German translations (rows 4-10). ---

$ws.Range("D4:D10").Font.Color = 0

# --- Column widths. ---

$ws.Columns.Item(1).ColumnWidth = 19.67
$ws.Columns.Item(3).ColumnWidth = 27.67
$ws.Columns.Item(4).ColumnWidth = 41

# --- Active selection, matching the authored file. ---

$ws.Range("C17").Select() | Out-Null
